$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) and Volume(1h) (column E) values for the
# refreshed symbol list snapshot. Values that look numeric (including
# percentages) are written with a leading apostrophe so Excel stores
# them as literal text, matching the original inline-string cells.
$ws.Range("D2").Value = "'309.70"
$ws.Range("E2").Value = "'0.21%"
$ws.Range("D3").Value = "'41.00"
$ws.Range("E3").Value = "'-1.11%"
$ws.Range("D4").Value = "'5.202"
$ws.Range("E4").Value = "'1.64%"
$ws.Range("E5").Value = "'0.19%"
$ws.Range("D6").Value = "'1.746"
$ws.Range("E6").Value = "'7.69%"
$ws.Range("D7").Value = "'0.9173"
$ws.Range("E7").Value = "'1.24%"
$ws.Range("D9").Value = "'0.1263"
$ws.Range("E9").Value = "'13.06%"
$ws.Range("D10").Value = "'0.1822"
$ws.Range("E10").Value = "'0.62%"
$ws.Range("D11").Value = "'0.09146"
$ws.Range("E11").Value = "'0.73%"
$ws.Range("D12").Value = "'0.04160"
$ws.Range("E12").Value = "'-1.90%"
$ws.Range("D13").Value = "'0.1051"
$ws.Range("E13").Value = "'-0.10%"
$ws.Range("D14").Value = "'0.001281"
$ws.Range("E14").Value = "'2.39%"
$ws.Range("D15").Value = "'0.005896"
$ws.Range("E15").Value = "'2.92%"
$ws.Range("E16").Value = "'0.26%"
$ws.Range("D17").Value = "'4.300"
$ws.Range("E17").Value = "'0.79%"
$ws.Range("D19").Value = "'7.437"
$ws.Range("E19").Value = "'11.75%"
$ws.Range("E20").Value = "'-0.71%"
$ws.Range("E21").Value = "'-0.37%"
$ws.Range("D22").Value = "'0.04015"
$ws.Range("E22").Value = "'-1.01%"
$ws.Range("D23").Value = "'0.001268"
$ws.Range("E23").Value = "'0.22%"
$ws.Range("D24").Value = "'0.004097"
$ws.Range("E24").Value = "'1.40%"
$ws.Range("E25").Value = "'0.26%"
$ws.Range("D38").Value = "'0.02521"
$ws.Range("E38").Value = "'3.66%"
$ws.Range("D39").Value = "'0.05300"
$ws.Range("E39").Value = "'0.71%"
$ws.Range("D40").Value = "'0.007857"
$ws.Range("E40").Value = "'0.81%"
$ws.Range("D41").Value = "'0.1310"
$ws.Range("E41").Value = "'0.64%"
$ws.Range("E42").Value = "'2.01%"
$ws.Range("E43").Value = "'-3.85%"
$ws.Range("D44").Value = "'0.008129"
$ws.Range("E44").Value = "'7.48%"
$ws.Range("E45").Value = "'-0.27%"
$ws.Range("D46").Value = "'0.00006792"
$ws.Range("E46").Value = "'0.21%"
$ws.Range("E47").Value = "'0.26%"
$ws.Range("D48").Value = "'0.2275"
$ws.Range("E48").Value = "'241.09%"
$ws.Range("E50").Value = "'0.26%"
$ws.Range("E51").Value = "'0.26%"
